$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list (Price / Volume(1h) columns) with the latest scrape.
# A leading apostrophe forces Excel to store the value as text instead of
# auto-converting number-like strings (e.g. "95.20" -> 95.2, "1.913.41"
# already non-numeric) so the stored text matches exactly.
$ws.Range("D2").Value = "30.141.62"
$ws.Range("E2").Value = "  +5.65%  "
$ws.Range("D3").Value = "1.914.52"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'329.55"
$ws.Range("E5").Value = "  +4.85%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "'0.5186"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("D8").Value = "'0.4061"
$ws.Range("E8").Value = "  +3.80%  "
$ws.Range("D9").Value = "'0.08487"
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").Value = "'1.124"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("D11").Value = "'42.75"
$ws.Range("D12").Value = "'23.15"
$ws.Range("E12").Value = "  +14.06%  "
$ws.Range("D13").Value = "'6.434"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").Value = "1.918.46"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").Value = "'7.383"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "'95.20"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "'0.00001114"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "'0.06695"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'18.41"
$ws.Range("E20").Value = "  +4.66%  "
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'6.005"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "30.149.82"
$ws.Range("E23").Value = "  +5.58%  "
$ws.Range("D24").Value = "'11.32"
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("D25").Value = "'2.237"
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("D26").Value = "2.135.96"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").Value = "'161.84"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").Value = "'21.13"
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("D29").Value = "'2.411"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").Value = "'128.79"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").Value = "'1.095"
$ws.Range("E31").Value = "  +5.18%  "
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").Value = "'6.019"
$ws.Range("E33").Value = "  +4.27%  "
$ws.Range("D34").Value = "'3.631"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'0.02489"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").Value = "'0.06577"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "'0.2209"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.227"
$ws.Range("E38").Value = "  +3.40%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.172"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").Value = "'11.87"
$ws.Range("E40").Value = "  +6.98%  "
$ws.Range("D41").Value = "'8.796"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").Value = "'0.6522"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'0.6136"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").Value = "'13.36"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").Value = "'3.743"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").Value = "'2.076"
$ws.Range("E47").Value = "  +3.65%  "
$ws.Range("D48").Value = "'1.244"
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").Value = "'124.14"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "'1.163"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").Value = "'79.44"
